# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# Price-looking numeric strings (e.g. "304.42") are written with a leading
# apostrophe so Excel keeps them as literal text (matching the sheet's
# existing inlineStr "Price" column) instead of silently coercing them to
# floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.238.20'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '1.605.14'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = '''304.42'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').Value = '''52.34'
$ws.Range('E8').Value = '  +4.85%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').Value = '''1.273'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '''22.92'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('D14').Value = '''6.600'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = '''7.386'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '''0.00001251'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '1.605.96'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '''93.88'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = '''0.06918'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('D24').Value = '23.242.01'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = '''3.103'
$ws.Range('E25').Value = '  +10.32%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '''2.450'
$ws.Range('E26').Value = '  +3.77%  '
$ws.Range('D27').Value = '''21.18'
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('D28').Value = '''149.97'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '''5.278'
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').Value = '''135.07'
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('D31').Value = '''2.394'
$ws.Range('E31').Value = '  +2.07%  '
$ws.Range('D32').Value = '''6.752'
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('D33').Value = '1.780.78'
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').Value = '''0.9642'
$ws.Range('E34').Value = '  +0.27%  '
$ws.Range('D35').Value = '''0.07489'
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('E36').Value = '  +0.74%  '
$ws.Range('D37').Value = '''0.02760'
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('D38').Value = '''0.2517'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = '''6.117'
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('D40').Value = '''0.08799'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('D41').Value = '''1.423'
$ws.Range('E41').Value = '  +4.24%  '
$ws.Range('D42').Value = '''0.7097'
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('D43').Value = '''12.48'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').Value = '''15.80'
$ws.Range('E44').Value = '  +3.73%  '
$ws.Range('D45').Value = '''0.6534'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').Value = '''2.331'
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = '''4.009'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').Value = '''133.91'
$ws.Range('D49').Value = '''0.07936'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').Value = '''1.209'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').Value = '''1.191'
$ws.Range('E51').Value = '  -3.04%  '
